$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.596.28'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.514.44'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.67'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.40'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.586'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.14%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.539'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.65'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0813'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.75'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.112'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.904.20'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.76'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +8.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.522.09'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.863'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.572.59'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.97'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0973'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.42%  '
$ws.Range("E21").Value = '  -1.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.55'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.27'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.65%  '
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.04'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.17'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.33'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +10.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.18'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.61'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.95'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '153.78'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("E33").Value = '  +5.06%  '
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0789'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.60%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.08'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.07%  '
$ws.Range("E37").Value = '  -4.89%  '
$ws.Range("E38").Value = '  -1.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.72'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.24%  '
$ws.Range("E40").Value = '  +0.57%  '
$ws.Range("E41").Value = '  +0.77%  '
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("E43").Value = '  -0.52%  '
$ws.Range("E44").Value = '  -0.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.037.24'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.67'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.98'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.758.04'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.78%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.190'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.88'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.05%  '

Write-Output "applied changes"
